$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$d5 = 'Objectives
id="Par1">Rapid and accurate diagnosis of coronavirus disease 2019 (COVID-19) is critical during the epidemic.

 We aim to identify differences in CT imaging and clinical manifestations between pneumonia patients with and without COVID-19, and to develop and validate a diagnostic model for COVID-19 based on radiological semantic and clinical features alone.


Methods
id="Par2">A consecutive cohort of 70 COVID-19 and 66 non-COVID-19 pneumonia patients were retrospectively recruited from five institutions.

 Patients were divided into primary (n = 98) and validation (n = 38) cohorts.

 The chi-square test, Student’s t test, and Kruskal-Wallis H test were performed, comparing 1745 lesions and 67 features in the two groups.

 Three models were constructed using radiological semantic and clinical features through multivariate logistic regression.

 Diagnostic efficacies of developed models were quantified by receiver operating characteristic curve.

 Clinical usage was evaluated by decision curve analysis and nomogram.


Results
id="Par3">Eighteen radiological semantic features and seventeen clinical features were identified to be significantly different.

 Besides ground-glass opacities (p = 0.032) and consolidation (p = 0.001) in the lung periphery, the lesion size (1–3 cm) is also significant for the diagnosis of COVID-19 (p = 0.027).

 Lung score presents no significant difference (p = 0.417).

 Three diagnostic models achieved an area under the curve value as high as 0.986 (95% CI 0.966~1.000).

 The clinical and radiological semantic models provided a better diagnostic performance and more considerable net benefits.


Conclusions
id="Par4">Based on CT imaging and clinical manifestations alone, the pneumonia patients with and without COVID-19 can be distinguished.

 A model composed of radiological semantic and clinical features has an excellent performance for the diagnosis of COVID-19.
Key Points
id="Par5">
            • Based on CT imaging and clinical manifestations alone, the pneumonia patients with and without COVID-19 can be distinguished.


          
Electronic supplementary material
The online version of this article (10.1007/s00330-020-06829-2) contains supplementary material, which is available to authorized users.


'
$e5 = '[Xiaofeng%Chen%NULL%2,      Yanyan%Tang%NULL%1,      Yongkang%Mo%NULL%1,      Shengkai%Li%NULL%1,      Daiying%Lin%NULL%1,      Zhijian%Yang%NULL%1,      Zhiqi%Yang%NULL%1,      Hongfu%Sun%NULL%1,      Jinming%Qiu%NULL%1,      Yuting%Liao%NULL%1,      Jianning%Xiao%NULL%1,      Xiangguang%Chen%NULL%1,      Xianheng%Wu%NULL%1,      Renhua%Wu%NULL%1,      Zhuozhi%Dai%zhuozhi@ualberta.ca%1,     Xiaofeng%Chen%NULL%1,     Yanyan%Tang%NULL%1,     Yongkang%Mo%NULL%1,     Shengkai%Li%NULL%1,     Daiying%Lin%NULL%1,     Zhijian%Yang%NULL%1,     Zhiqi%Yang%NULL%1,     Hongfu%Sun%NULL%1,     Jinming%Qiu%NULL%1,     Yuting%Liao%NULL%1,     Jianning%Xiao%NULL%1,     Xiangguang%Chen%NULL%1,     Xianheng%Wu%NULL%1,     Renhua%Wu%NULL%1,     Zhuozhi%Dai%zhuozhi@ualberta.ca%1]'
$d8 = 'Purpose
id="Par1">To evaluate the diagnostic performance of chest CT to differentiate coronavirus disease 2019 (COVID-19) pneumonia in non-high-epidemic area in Japan.


Materials and methods
id="Par2">This retrospective study included 21 patients clinically suspected COVID-19 pneumonia and underwent chest CT more than 3 days after the symptom onset: six patients confirmed COVID-19 pneumonia by real-time reverse-transcription polymerase chain reaction (RT-PCR) and 15 patients proved uninfected.

 Using a Likert scale and its receiver operating characteristic curve analysis, two radiologists (R1/R2) evaluated the diagnostic performance of the five CT criteria: (1) ground glass opacity (GGO)-predominant lesions, (2) GGO- and peripheral-predominant lesions, (3) bilateral GGO-predominant lesions; (4) bilateral GGO- and peripheral-predominant lesions, and (5) bilateral GGO- and peripheral-predominant lesions without nodules, airway abnormalities, pleural effusion, and mediastinal lymphadenopathy.


Results
id="Par3">All patients confirmed COVID-19 pneumonia had bilateral GGO- and peripheral-predominant lesions without airway abnormalities, mediastinal lymphadenopathy, and pleural effusion.

 The five CT criteria showed moderate to excellent diagnostic performance with area under the curves (AUCs) ranging 0.77–0.88 for R1 and 0.78–0.92 for R2. The criterion (e) showed the highest AUC.


Conclusion
id="Par4">Chest CT would play a supplemental role to differentiate COVID-19 pneumonia from other respiratory diseases presenting with similar symptoms in a clinical setting.


Electronic supplementary material
The online version of this article (10.1007/s11604-020-00958-w) contains supplementary material, which is available to authorized users.


'
$e8 = '[Yuki%Himoto%yukihimoto@gmail.com%2,      Akihiko%Sakata%NULL%1,      Mitsuhiro%Kirita%NULL%1,      Takashi%Hiroi%NULL%1,      Ken-ichiro%Kobayashi%NULL%1,      Kenji%Kubo%NULL%1,      Hyunjin%Kim%NULL%1,      Azusa%Nishimoto%NULL%1,      Chikara%Maeda%NULL%1,      Akira%Kawamura%NULL%1,      Nobuhiro%Komiya%NULL%1,      Shigeaki%Umeoka%NULL%1,     Yuki%Himoto%yukihimoto@gmail.com%1,     Akihiko%Sakata%NULL%1,     Mitsuhiro%Kirita%NULL%1,     Takashi%Hiroi%NULL%1,     Ken-ichiro%Kobayashi%NULL%1,     Kenji%Kubo%NULL%1,     Hyunjin%Kim%NULL%1,     Azusa%Nishimoto%NULL%1,     Chikara%Maeda%NULL%1,     Akira%Kawamura%NULL%1,     Nobuhiro%Komiya%NULL%1,     Shigeaki%Umeoka%NULL%1]'

$ws.Range("D5").Value = $d5
$ws.Range("E5").Value = $e5
$ws.Range("D8").Value = $d8
$ws.Range("E8").Value = $e8
